$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new data row (Verdict "Neutral" set first so it gets the next
#     shared-string slot before "Down", matching the target shared string order) ---
$ws.Range("A4").Value = 42633.886388888888
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Neutral"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Random"
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.49
$ws.Range("S4").NumberFormat = "0.00%"
$ws.Range("S4").Value = 0.088800000000000004
$ws.Range("T4").Value = -1.1000000000000001
$ws.Range("U4").Value = 2.2999999999999998
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0

# --- Row 3: two new trailing columns (sentiment delta + up/down verdict) ---
$ws.Range("X3").Value = -0.19999900000000181
$ws.Range("Y3").Value = "Down"

# --- Column C ("Verdict") widened slightly now that "Neutral" is present ---
# (COM ColumnWidth only accepts pixel-quantized widths; 6.8 lands on the grid
# point closest to the real bestFit width of 7.7109375 that Excel would compute.)
$ws.Columns("C").ColumnWidth = 6.8
